$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45044
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 833

$ws.Range("D3").Value = 45043
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("S3").Value = 833

$ws.Range("D4").Value = 45030
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("S4").Value = 1000

$ws.Range("D5").Value = 45041
$ws.Range("M5").Value = 60

$ws.Range("D6").Value = 45001
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 17000
$ws.Range("P6").Value = 17500
$ws.Range("S6").Value = 972

$ws.Range("D7").Value = 45028
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("S7").Value = 1000

$ws.Range("D8").Value = 45014
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("S8").Value = 1000

$ws.Range("D9").Value = 45050
$ws.Range("M9").Value = 40
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 14000
$ws.Range("P9").Value = 14000
$ws.Range("R9").Value = 'Región Metropolitana'
$ws.Range("S9").Value = 778

$ws.Range("D10").Value = 45002
$ws.Range("M10").Value = 30

$ws.Range("D11").Value = 45049
$ws.Range("M11").Value = 80

$ws.Range("D12").Value = 44999
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 17500
$ws.Range("S12").Value = 972

$ws.Range("D13").Value = 45020
$ws.Range("M13").Value = 50
$ws.Range("Q13").Value = '$/caja 16 kilos'
$ws.Range("R13").Value = 'Provincia de Los Andes'
$ws.Range("S13").Value = 938
$ws.Range("T13").Value = 16

$ws.Range("D14").Value = 45062
$ws.Range("M14").Value = 90
$ws.Range("N14").Value = 13000
$ws.Range("O14").Value = 14000
$ws.Range("P14").Value = 13444
$ws.Range("S14").Value = 747

$ws.Range("D15").Value = 45033
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 15500
$ws.Range("S15").Value = 861

$ws.Range("D16").Value = 45036
$ws.Range("M16").Value = 60
$ws.Range("O16").Value = 16000
$ws.Range("P16").Value = 15500
$ws.Range("Q16").Value = '$/caja 18 kilos'
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("S16").Value = 861
$ws.Range("T16").Value = 18

$ws.Range("D17").Value = 45021
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 15500
$ws.Range("R17").Value = 'Provincia de Los Andes'
$ws.Range("S17").Value = 861

$ws.Range("D18").Value = 45037
$ws.Range("N18").Value = 16000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 16000
$ws.Range("S18").Value = 889
